$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final parent "Post Treatment" measurements (column D) for rows 2-14
$ws.Range("D2").Value = 44
$ws.Range("D3").Value = 45
$ws.Range("D4").Value = 44
$ws.Range("D5").Value = 44
$ws.Range("D6").Value = 60
$ws.Range("D7").Value = 52
$ws.Range("D8").Value = 49
$ws.Range("D9").Value = 67
$ws.Range("D10").Value = 53
$ws.Range("D11").Value = 44
$ws.Range("D12").Value = 67
$ws.Range("D13").Value = 58
$ws.Range("D14").Value = 54

# Move selection to reflect where data entry ended (one row below the last entry)
$ws.Range("D15").Select()
